$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(30).Insert()

$ws.AutoFilterMode = $false
$ws.Range("A4:CD4").AutoFilter()

$n = $wb.Names.Item(1)
$n.RefersTo = "=Sheet1!`$A`$4:`$CD`$4"

$ws.Rows.Item(30).Insert()
Write-Host "UsedRange after row insert:" $ws.UsedRange.Address()
Write-Host "A30:" $ws.Range("A30").Value()
Write-Host "A31:" $ws.Range("A31").Value()
